# Release 02 of "The Adventures of an Adventurer" + Release and Planing Docs
#
# Product backlog updates:
#  - Row 13 ("Design new Enemies") is re-scoped to the second level, its
#    Vertics user-story text is reworded, and its SP estimate drops 8 -> 5.
#  - A batch of already-delivered level-1 backlog items (rows 10,14,16,18-25)
#    get a "-" marker in column H, and row 15 gets marked "ok" (Finished).
#  - A brand new backlog item "Extended Tutorial Bossfight" is appended as
#    row 49.
#  - The sheet's saved view/selection moves to reflect where the author was
#    last working (top-left C7, active cell F14).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Row 13: re-scope "Design new Enemies" to the second level ---
$ws.Range("C13").Value = "Design new Enemies for second level"
$ws.Range("E13").Value = "As a Player I want to fight against new enemies in the second level. That means 10 new enemies per level. I also want to see some familiar enemies which I already met in previous games."
$ws.Range("F13").Value = 5

# --- Mark already finished level-1 items with "-" in column H ---
$ws.Range("H10").Value = "-"
$ws.Range("H14").Value = "-"
$ws.Range("H16").Value = "-"
$ws.Range("H18").Value = "-"
$ws.Range("H19").Value = "-"
$ws.Range("H20").Value = "-"
$ws.Range("H21").Value = "-"
$ws.Range("H22").Value = "-"
$ws.Range("H23").Value = "-"
$ws.Range("H24").Value = "-"
$ws.Range("H25").Value = "-"

# Row 15 ("Convert Level form Paper to Unity") is now finished
$ws.Range("G15").Value = "ok"

# --- New backlog item: Extended Tutorial Bossfight (row 49) ---
$ws.Range("C49").Value = "Extended Tutorial Bossfight"
$ws.Range("D49").Value = "!!"
$ws.Range("E49").Value = "As a player I want the first boss I fight against to be able to use multiple different attacks instead of just one."
$ws.Range("F49").Value = 5
$ws.Range("G49").Value = "ok"

# --- Restore the author's last-viewed window position/selection ---
$ws.Range("F14").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 3
